$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename SKU codes in column A (rows 2-101) from GSQL00001..GSQL00100
# to TWNT00001..TWNT00100, preserving the numeric suffix.
for ($i = 1; $i -le 100; $i++) {
    $num = "{0:D5}" -f $i
    $ws.Cells.Item($i + 1, 1).Value = "TWNT$num"
}
